# Updates the "Count down" control-flow code sample that appears on every
# slide of this deck, renaming the C#-style API calls used in the sample
# code to their snake_case equivalents:
#   Write(    -> write(
#   ToInt32   -> to_integer
#   ReadLine  -> read_line
#   WriteLine -> write_line
#
# The code sample lives in a text box named "TextBox 15" inside a group
# named "Group 11" on every slide (1-10). Each renamed token is already
# its own complete run in the OOXML, so each replacement is done as a
# whole-run Characters() replace (matching the exact existing run text)
# rather than touching a sub-string of a run - that way PowerPoint does
# not fragment the run into extra pieces that aren't in the target.
#
# Two wrinkles:
#  * "for (max = ToInt32(" is a single fused run on slides 1/4/7/10, but
#    split into two runs - "for (" and "max = ToInt32(" - on slides
#    2/3/5/6/8/9 (the two runs render identically but carry different
#    direct formatting). Both forms flatten to the same plain text, so
#    which whole-run pattern to search for is chosen per slide.
#  * "    WriteLine(max);" is a single fused run on every slide except
#    slide 10, where it is split into "    WriteLine" + "(max);".
#  * The COM TextRange.Text getter normalises the smart/curly close
#    quote (U+201D) used in `WriteLine(”Blast off !");` down to a plain
#    ASCII quote when read back, even though the underlying run still
#    stores the curly quote. So the search pattern uses the ASCII quote
#    while the replacement text reconstructs the original curly quote.

function Update-CodeSnippet($tr, $forFused, $writeLineFused) {
    $quote = [string][char]0x201d
    $dquote = [string][char]0x22

    # "WriteLine(<curly-quote>Blast off !");" -> "write_line(<curly-quote>Blast off !");"
    $pat = "WriteLine(" + $dquote + "Blast off !" + $dquote + ");"
    $t = $tr.Text
    $idx = $t.IndexOf($pat)
    if ($idx -ge 0) {
        $new = "write_line(" + $quote + "Blast off !" + $dquote + ");"
        $tr.Characters($idx + 1, $pat.Length).Text = $new
    }

    # "    WriteLine(max);" loop body - fused or split variant
    if ($writeLineFused) {
        $pat = "    WriteLine(max);"
        $t = $tr.Text
        $idx = $t.IndexOf($pat)
        if ($idx -ge 0) {
            $tr.Characters($idx + 1, $pat.Length).Text = "    write_line(max);"
        }
    } else {
        $pat = "    WriteLine"
        $t = $tr.Text
        $idx = $t.IndexOf($pat)
        if ($idx -ge 0) {
            $tr.Characters($idx + 1, $pat.Length).Text = "    write_line"
        }
    }

    # "for (max = ToInt32(" fused, or "max = ToInt32(" on its own (with a
    # separate preceding "for (" run left untouched)
    if ($forFused) {
        $pat = "for (max = ToInt32("
        $t = $tr.Text
        $idx = $t.IndexOf($pat)
        if ($idx -ge 0) {
            $tr.Characters($idx + 1, $pat.Length).Text = "for (max = to_integer("
        }
    } else {
        $pat = "max = ToInt32("
        $t = $tr.Text
        $idx = $t.IndexOf($pat)
        if ($idx -ge 0) {
            $tr.Characters($idx + 1, $pat.Length).Text = "max = to_integer("
        }
    }

    # "ReadLine" (own run on every slide)
    $pat = "ReadLine"
    $t = $tr.Text
    $idx = $t.IndexOf($pat)
    if ($idx -ge 0) {
        $tr.Characters($idx + 1, $pat.Length).Text = "read_line"
    }

    # 'Write("Count down from: ");' (own run on every slide)
    $pat = "Write(" + $dquote + "Count down from: " + $dquote + ");"
    $t = $tr.Text
    $idx = $t.IndexOf($pat)
    if ($idx -ge 0) {
        $new = "write(" + $dquote + "Count down from: " + $dquote + ");"
        $tr.Characters($idx + 1, $pat.Length).Text = $new
    }
}

$p = $ppt.ActivePresentation

# Slides whose "for (... ToInt32(" run is fused into a single run
# (vs. split into a separate "for (" run + "max = ToInt32(" run).
$forFusedSlides = @(1, 4, 7, 10)
# Slides whose "    WriteLine(max);" loop-body run is fused into a
# single run (vs. split into "    WriteLine" + "(max);").
$writeLineSplitSlides = @(10)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)
        if ($shp.Name -eq "Group 11") {
            $tb = $shp.GroupItems.Item("TextBox 15")
            $forFused = $forFusedSlides -contains $i
            $writeLineFused = -not ($writeLineSplitSlides -contains $i)
            Update-CodeSnippet $tb.TextFrame.TextRange $forFused $writeLineFused
        }
    }
}
